$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 10:05"

# Row 5 - Rusia
$ws.Range("B5").Value = 299941
$ws.Range("C5").Value = 9263
$ws.Range("D5").Value = 76130
$ws.Range("E5").Value = 220974
$ws.Range("G5").Value = 115
$ws.Range("H5").Value = 2837

# Row 34 - Polonia
$ws.Range("B34").Value = 19080
$ws.Range("C34").Value = 195
$ws.Range("D34").Value = 7903
$ws.Range("E34").Value = 10236
$ws.Range("G34").Value = 5
$ws.Range("H34").Value = 941

# Row 35 - Ucrania
$ws.Range("B35").Value = 18876
$ws.Range("C35").Value = 260
$ws.Range("D35").Value = 5632
$ws.Range("E35").Value = 12696
$ws.Range("G35").Value = 13
$ws.Range("H35").Value = 548

# Row 88 - Estonia
$ws.Range("B88").Value = 1791
$ws.Range("C88").Value = 7
$ws.Range("E88").Value = 789

# Row 95 - Eslovaquia
$ws.Range("D95").Value = 1192
$ws.Range("E95").Value = 275
